# Word COM-interop script: "actualizacion de perfil de colaboradores"
#
# This applies two related changes to the collaborators table:
#  1. The name "Jhono" + "rys" (two separate runs split by a stray
#     "_GoBack" bookmark) is merged into a single run "Jhonorys",
#     removing the old bookmark in the process.
#  2. A "_GoBack" bookmark is (re)created around the run holding
#     "Especialista en Información y Estadística " in the role column,
#     reflecting that this is where the document was last edited.

$d = $word.ActiveDocument

# --- Step 1: merge "Jhono" + "rys" into "Jhonorys" -------------------
# Word's Find searches the logical text stream, so "Jhono" + bookmark +
# "rys" already reads as "Jhonorys". Doing a same-text Find & Replace
# over that span collapses it into a single run and drops the
# now-redundant bookmark that used to sit between the two runs.
$d.Content.Find.Execute("Jhonorys", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Jhonorys", 2) | Out-Null

# --- Step 2: move the "_GoBack" bookmark onto the new edit location --
# Locate the run for "Especialista en Información y Estadística "
# (including its trailing space) and wrap it with a "_GoBack" bookmark,
# matching where Word would leave its last-edit marker.
$r = $d.Content
$r.Find.Execute("Especialista en Información y Estadística ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
